$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a new row for "Disclosure risk" so the (alphabetically sorted)
#        term table stays sorted: it belongs right before "Experimental data",
#        currently on row 22.
$ws.Rows("22:22").Insert()

$ws.Range("A22").Value = "Disclosure risk"
$ws.Range("C22").Value = "The risk of re-identifying a participant and the harm that may come from that disclosure."

# --- 2) Update the "Indirect identifiers" row (now shifted down to row 30):
#        add the "quasi-identifiers" alternate term and tweak the definition
#        wording (drop "or if category numbers are small").
$ws.Range("B30").Value = "quasi-identifiers"
$ws.Range("C30").Value = "These variables do not alone identify a particular individual (e.g., ethnicity, gender), but if combined with other information, they could be used to identify a participant"

# --- 3) Update the "Unique participant identifier" row (now shifted down to
#        row 63): append ", proxy ID" to the list of other terms.
$ws.Range("B63").Value = "study ID, site ID, unique identifier (UID), subject ID, participant code, record id, proxy ID"

# --- 4) Refresh the sheet's persisted sort state/range so it covers the new
#        row (A2:C65 -> A2:C66) just like Excel does after a Data > Sort.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A66"))
$sortObj.SetRange($ws.Range("A2:C66"))
$sortObj.Header = 0
$sortObj.Apply()

# --- 5) Restore the selection Excel leaves behind after this kind of edit.
$ws.Range("C13").Select()
